# 苏州-漫展信息.xlsx — apply "gh-pages output generated at 456a3b4" update.
#
# 昆山·ETHEREAL动漫游戏展（免费展） was cancelled: its title now carries a
# （取消） suffix and its "最低票价" (min. ticket price) cell switches from a
# numeric price to the text "不可售" (not for sale). That row/event lives on
# both the "展览" sheet and the combined "全部类型" sheet (same row numbers
# in this workbook), so the edit is applied identically on both.
#
# The remaining edits are routine "想去人数" (want-to-go count) refreshes —
# small incremental bumps to column F on several rows, spread across the
# "展览", "演出" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

function Set-WantCount {
    param($ws, $row, $value)
    $ws.Range("F$row").Value = $value
}

# ---- 展览 (Exhibitions) ----
$wsExpo = $wb.Worksheets.Item("展览")

$wsExpo.Range("C2").Value = "昆山·ETHEREAL动漫游戏展（免费展）（取消）"
$wsExpo.Range("G2").Value = "不可售"

Set-WantCount $wsExpo 7 14678
Set-WantCount $wsExpo 9 670
Set-WantCount $wsExpo 10 15165
Set-WantCount $wsExpo 12 8605
Set-WantCount $wsExpo 13 317
Set-WantCount $wsExpo 16 177
Set-WantCount $wsExpo 24 1071
Set-WantCount $wsExpo 26 9
Set-WantCount $wsExpo 30 415
Set-WantCount $wsExpo 31 22
Set-WantCount $wsExpo 33 229
Set-WantCount $wsExpo 35 417
Set-WantCount $wsExpo 37 5314

# ---- 演出 (Performances) ----
$wsShow = $wb.Worksheets.Item("演出")

Set-WantCount $wsShow 3 58

# ---- 全部类型 (All types, combined) ----
$wsAll = $wb.Worksheets.Item("全部类型")

$wsAll.Range("C2").Value = "昆山·ETHEREAL动漫游戏展（免费展）（取消）"
$wsAll.Range("G2").Value = "不可售"

Set-WantCount $wsAll 7 14678
Set-WantCount $wsAll 9 670
Set-WantCount $wsAll 10 15165
Set-WantCount $wsAll 12 8605
Set-WantCount $wsAll 13 317
Set-WantCount $wsAll 17 177
Set-WantCount $wsAll 25 1071
Set-WantCount $wsAll 27 9
Set-WantCount $wsAll 31 58
Set-WantCount $wsAll 33 415
Set-WantCount $wsAll 34 22
Set-WantCount $wsAll 36 229
Set-WantCount $wsAll 38 417
Set-WantCount $wsAll 40 5314
